$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(3,4,5,6,7,8,10,13,14,15)

$data = @(
    @(0.04960438181704774, 0.1922055140664654, 0.1321776307253534, 1.091459290066467, 0.6332093872325828, 0.6766019342329201, 0.1242006871250769, 1.062633181663088, 1.813524949418252, 2.605098433374792),
    @(0.04402263596996647, 0.1941193837986361, 0.1343682001590979, 1.0744093040189, 0.6086734561348379, 0.6698847976916227, 0.1268758216931052, 0.9554597755403762, 1.665631701753284, 2.53862550547575),
    @(0.04061103459267201, 0.1953663295942807, 0.135788424637275, 1.064696601847416, 0.5940811851556731, 0.6661616613001797, 0.1286108358439733, 0.8894764630605465, 1.574993815108371, 2.49963818291252),
    @(0.03922465822222421, 0.1958925342865392, 0.1363860647016674, 1.060928236255336, 0.5882531194770593, 0.6647452239511438, 0.1293410668729664, 0.8625454717536769, 1.538105106546169, 2.484208693179198),
    @(0.03899468438234521, 0.1959810013215311, 0.1364864424054107, 1.060313941860457, 0.58729251003831, 0.6645161085388338, 0.1294637207478431, 0.8580711259283333, 1.53198273238246, 2.4816742710874),
    @(0.04059232179156425, 0.1953733530280513, 0.1357964081795836, 1.06464501314484, 0.5940021070887695, 0.666142150872119, 0.1286205901365527, 0.8891134299239809, 1.574496125015344, 2.49942824211962),
    @(0.04767652531985789, 0.1928505122191009, 0.1329173158402202, 1.085423238861466, 0.624650934447132, 0.67420251903701, 0.1251038491690748, 1.025718346478854, 1.762498840636766, 2.581798474916212),
    @(0.06169532997853366, 0.1884728724651286, 0.1278691350506883, 1.13219351241392, 0.6885321986151212, 0.6931994145703868, 0.118943537577096, 1.292083054324664, 2.132339682909901, 2.75789747893117),
    @(0.07207708813358238, 0.1856034903473098, 0.1245258330929258, 1.170268475281375, 0.7378122968161165, 0.7091139055110034, 0.1148692250890404, 1.486738951118824, 2.404569657386105, 2.896277753092363),
    @(0.07681896428178447, 0.184373366885465, 0.1230845211492422, 1.188404996746243, 0.7607501605692732, 0.7167817273580113, 0.1131144263256338, 1.575043111364508, 2.528483716617302, 2.961211224013482),
    @(0.07861741114410847, 0.1839183572319172, 0.1225502019663467, 1.195390792268327, 0.7695115737348601, 0.7197470944199438, 0.1124641701217666, 1.608443892844605, 2.575413648292511, 2.986086850347931),
    @(0.07822995811930866, 0.184015870835081, 0.1226647663281411, 1.193881024402856, 0.7676212869877759, 0.7191057017497826, 0.1126035796976437, 1.601252174394858, 2.565306235060973, 2.980716662915995),
    @(0.0769668671848649, 0.1843357162807333, 0.1230403322420981, 1.188977356564251, 0.7614694541710776, 0.7170244518857487, 0.1130606434282013, 1.577791792745927, 2.532344568083374, 2.963252001699857),
    @(0.07619355368045433, 0.1845330388656254, 0.1232718724581241, 1.185989087162994, 0.7577111070974638, 0.7157576696735646, 0.1133424656237816, 1.56341659031645, 2.512155293071089, 2.952591777521377),
    @(0.07176758488537871, 0.1856853946176162, 0.1246216301514973, 1.169099678053641, 0.7363237663994653, 0.7086214253071432, 0.114985895536778, 1.480962885793559, 2.396472714629397, 2.89207424041814),
    @(0.06905733692263993, 0.1864115822654888, 0.1254700659613417, 1.158947964258402, 0.7233369283091804, 0.7043533586557658, 0.1160193937765879, 1.430315299500478, 2.325521290693189, 2.855457796903238),
    @(0.06750027490159027, 0.1868363417168641, 0.1259655516619438, 1.153185742924606, 0.7159161637169404, 0.701938778588044, 0.1166231138710287, 1.401161176035288, 2.284719256762003, 2.834583602412181),
    @(0.06697338836849553, 0.1869813731481198, 0.1261345996848282, 1.151247923577088, 0.7134120057850168, 0.7011281601434689, 0.1168291148458547, 1.391286229317117, 2.270905788017217, 2.827547964815494),
    @(0.06934566072521875, 0.186333546049303, 0.1253789733241595, 1.160020681558706, 0.7247143322368288, 0.7048035297646891, 0.1159084152823411, 1.435709218909807, 2.333073462815321, 2.859336355031019),
    @(0.07733779122060014, 0.1842414764885731, 0.1229297077635598, 1.190414478869528, 0.7632743474262043, 0.7176340889106143, 0.1129260055200527, 1.58468373099096, 2.542026071143141, 2.968374002357393),
    @(0.08257745961650187, 0.1829372060391687, 0.1213958674025712, 1.210965980376528, 0.788914954164909, 0.7263794769361311, 0.1110599072491771, 1.681823922277346, 2.678623759850495, 3.041308386852393),
    @(0.07977943895765804, 0.1836275530203118, 0.1222083752983811, 1.199934179997967, 0.7751896912211862, 0.7216789205386362, 0.1120482555344404, 1.629999704804064, 2.605717251639078, 3.00222845273214),
    @(0.06921530612642357, 0.1863688035951654, 0.1254201322951725, 1.15953547517924, 0.7240914661032036, 0.7045998853408548, 0.1159585589210437, 1.433270740340461, 2.329659160076858, 2.857582307385712),
    @(0.05788880027417065, 0.189596197388342, 0.1291706874114815, 1.118891818592445, 0.6708416405862749, 0.6877174018005974, 0.1133424656237816, 1.220198960981151, 2.032184207337934, 2.980716662915995)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $rowVals[$j]
    }
}
